$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 8
# Leading apostrophe forces literal/text storage (matches the existing rows,
# which store even numeric-looking values like "22" as text).
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 2).Value = "احمد"
$ws.Cells.Item($row, 3).Value = "'22"
$ws.Cells.Item($row, 4).Value = "الجزائري"
$ws.Cells.Item($row, 5).Value = "الرحلة 1"
$ws.Cells.Item($row, 6).Value = "C3"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٤:٤٠:٠٤ م"
